# Adds the required "Experimental" boolean value (stored as text "true")
# to the Metadata sheet, and refreshes the "Date" value, as described by:
#   "added required experimental boolean element to valuesets"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = "Experimental" property; its value cell (B7) was empty and must
# become the text "true". Assigning the literal string "true"/"false"
# directly to .Value/.Value2 gets auto-coerced to an Excel Boolean, so we
# instead write a formula that evaluates to the text "true" and then
# collapse it down to a plain value in-place (Copy + PasteSpecial values),
# which keeps it a normal text cell using the same existing cell style.
$b7 = $ws.Cells.Item(7, 2)
$b7.Formula = "=""true"""
$b7.Copy()
$b7.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Row 8 = "Date" property; refresh its value to the new publication date.
$ws.Cells.Item(8, 2).Value = "2023-02-01T09:05:11-06:00"
